# Applies scheduled-runner market-data updates to the Moogle_Profits workbook sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2281.7334
$ws.Range("I33").Value = 562.8182
$ws.Range("K33").Value = 562.8182
$ws.Range("M33").Value = -333.8182

$ws.Range("H64").Value = 14184.143
$ws.Range("I64").Value = 9994.333000000001
$ws.Range("K64").Value = 9994.333000000001
$ws.Range("M64").Value = -9746.333000000001

$ws.Range("H67").Value = 14184.143
$ws.Range("I67").Value = 9994.333000000001
$ws.Range("K67").Value = 9994.333000000001
$ws.Range("M67").Value = -9136.333000000001

$ws.Range("H100").Value = 3999.5
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459

$ws.Range("H121").Value = 3790.2222
$ws.Range("J121").Value = 3790.2222
$ws.Range("L121").Value = 11370.6666
$ws.Range("N121").Value = -14864.6666

$ws.Range("H137").Value = 1929.1
$ws.Range("I137").Value = 1951.9412
$ws.Range("J137").Value = 1799.6666
$ws.Range("K137").Value = 5855.8236
$ws.Range("L137").Value = 5398.9998
$ws.Range("M137").Value = -3305.8236
$ws.Range("N137").Value = -10498.9998

$ws.Range("H141").Value = 1730.0834
$ws.Range("I141").Value = 1617.762
$ws.Range("J141").Value = 2516.3333
$ws.Range("K141").Value = 4853.286
$ws.Range("L141").Value = 7548.999899999999
$ws.Range("M141").Value = 326.7139999999999
$ws.Range("N141").Value = -17908.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10900.924
$ws.Range("I32").Value = 6427.92
$ws.Range("K32").Value = 6427.92
$ws.Range("M32").Value = -6140.92

$ws.Range("H110").Value = 1875
$ws.Range("I110").Value = 1607.8125
$ws.Range("K110").Value = 1607.8125
$ws.Range("M110").Value = 437.1875

$ws.Range("H122").Value = 3763.05
$ws.Range("I122").Value = 2921.543
$ws.Range("J122").Value = 9653.6
$ws.Range("K122").Value = 8764.629000000001
$ws.Range("L122").Value = 28960.8
$ws.Range("M122").Value = -6314.629000000001
$ws.Range("N122").Value = -33860.8

$ws.Range("H132").Value = 4303.864
$ws.Range("I132").Value = 2934.25
$ws.Range("K132").Value = 8802.75
$ws.Range("M132").Value = -6272.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1993.0714
$ws.Range("I94").Value = 834.6667
$ws.Range("J94").Value = 2861.875
$ws.Range("K94").Value = 834.6667
$ws.Range("L94").Value = 2861.875
$ws.Range("M94").Value = -383.6667
$ws.Range("N94").Value = -3763.875

$ws.Range("H105").Value = 980071
$ws.Range("I105").Value = 1605198.9
$ws.Range("J105").Value = 3308.75
$ws.Range("K105").Value = 1605198.9
$ws.Range("L105").Value = 3308.75
$ws.Range("M105").Value = -1603451.9
$ws.Range("N105").Value = -6802.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8356.972
$ws.Range("I31").Value = 3633.2593
$ws.Range("J31").Value = 24299.5
$ws.Range("K31").Value = 3633.2593
$ws.Range("L31").Value = 24299.5
$ws.Range("M31").Value = -3338.2593
$ws.Range("N31").Value = -24889.5

$ws.Range("H34").Value = 8356.972
$ws.Range("I34").Value = 3633.2593
$ws.Range("J34").Value = 24299.5
$ws.Range("K34").Value = 3633.2593
$ws.Range("L34").Value = 24299.5
$ws.Range("M34").Value = -3431.2593
$ws.Range("N34").Value = -24703.5

$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws.Range("H132").Value = 2804.5894
$ws.Range("I132").Value = 2616.5
$ws.Range("K132").Value = 7849.5
$ws.Range("M132").Value = -5319.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3025.2727
$ws.Range("I46").Value = 321.5
$ws.Range("K46").Value = 964.5
$ws.Range("M46").Value = -873.5

$ws.Range("H121").Value = 7074964
$ws.Range("I121").Value = 2066.3333
$ws.Range("J121").Value = 10611413
$ws.Range("K121").Value = 6198.999899999999
$ws.Range("L121").Value = 31834239
$ws.Range("M121").Value = -4888.999899999999
$ws.Range("N121").Value = -31836859

$ws.Range("H131").Value = 4694.08
$ws.Range("J131").Value = 5722
$ws.Range("L131").Value = 17166
$ws.Range("N131").Value = -27246

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 326.33334
$ws.Range("I107").Value = 263.42856
$ws.Range("J107").Value = 546.5
$ws.Range("K107").Value = 263.42856
$ws.Range("L107").Value = 546.5
$ws.Range("M107").Value = 1656.57144
$ws.Range("N107").Value = -4386.5

$ws.Range("H122").Value = 3449.75
$ws.Range("I122").Value = 2080.3684
$ws.Range("J122").Value = 7166.643
$ws.Range("K122").Value = 6241.1052
$ws.Range("L122").Value = 21499.929
$ws.Range("M122").Value = -3791.1052
$ws.Range("N122").Value = -26399.929

$ws.Range("H126").Value = 8235.137000000001
$ws.Range("I126").Value = 7945.467
$ws.Range("K126").Value = 23836.401
$ws.Range("M126").Value = -21366.401

$ws.Range("H132").Value = 6215.5625
$ws.Range("I132").Value = 4388.2144
$ws.Range("J132").Value = 19007
$ws.Range("K132").Value = 13164.6432
$ws.Range("L132").Value = 57021
$ws.Range("M132").Value = -10634.6432
$ws.Range("N132").Value = -62081

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3979.9
$ws.Range("I40").Value = 2588.8333
$ws.Range("K40").Value = 2588.8333
$ws.Range("M40").Value = -2452.8333

$ws.Range("H61").Value = 6357.6523
$ws.Range("I61").Value = 2414.4666
$ws.Range("J61").Value = 13751.125
$ws.Range("K61").Value = 2414.4666
$ws.Range("L61").Value = 13751.125
$ws.Range("M61").Value = -2212.4666
$ws.Range("N61").Value = -14155.125

$ws.Range("H68").Value = 7470.6523
$ws.Range("I68").Value = 4483.364
$ws.Range("K68").Value = 4483.364
$ws.Range("M68").Value = -3734.364

$ws.Range("H71").Value = 7470.6523
$ws.Range("I71").Value = 4483.364
$ws.Range("K71").Value = 22416.82
$ws.Range("M71").Value = -18672.82

$ws.Range("H113").Value = 6357.6523
$ws.Range("I113").Value = 2414.4666
$ws.Range("J113").Value = 13751.125
$ws.Range("K113").Value = 2414.4666
$ws.Range("L113").Value = 13751.125
$ws.Range("M113").Value = -244.4666000000002
$ws.Range("N113").Value = -18091.125

$ws.Range("H136").Value = 6503.5137
$ws.Range("I136").Value = 4432.6875
$ws.Range("J136").Value = 10326.577
$ws.Range("K136").Value = 13298.0625
$ws.Range("L136").Value = 30979.731
$ws.Range("M136").Value = -10748.0625
$ws.Range("N136").Value = -36079.731

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 98354.75
$ws.Range("J46").Value = 98354.75
$ws.Range("L46").Value = 98354.75
$ws.Range("N46").Value = -98816.75

$ws.Range("H100").Value = 686.51514
$ws.Range("I100").Value = 445.25
$ws.Range("K100").Value = 890.5
$ws.Range("M100").Value = -349.5

$ws.Range("H107").Value = 3080
$ws.Range("I107").Value = 2901.9666
$ws.Range("J107").Value = 3413.8125
$ws.Range("K107").Value = 8705.899800000001
$ws.Range("L107").Value = 10241.4375
$ws.Range("M107").Value = -6785.899800000001
$ws.Range("N107").Value = -14081.4375

$ws.Range("H126").Value = 1844.9166
$ws.Range("I126").Value = 1113.9
$ws.Range("K126").Value = 3341.7
$ws.Range("M126").Value = -871.7000000000003

$ws.Range("H132").Value = 5287.0527
$ws.Range("I132").Value = 4064.7576
$ws.Range("K132").Value = 12194.2728
$ws.Range("M132").Value = -9664.272799999999

$ws.Range("H134").Value = 98354.75
$ws.Range("J134").Value = 98354.75
$ws.Range("L134").Value = 295064.25
$ws.Range("N134").Value = -300134.25

$ws.Range("H136").Value = 6493.0835
$ws.Range("I136").Value = 4895.4644
$ws.Range("J136").Value = 12084.75
$ws.Range("K136").Value = 14686.3932
$ws.Range("L136").Value = 36254.25
$ws.Range("M136").Value = -12136.3932
$ws.Range("N136").Value = -41354.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("N29").ClearContents()
